# Applies the "Add files via upload" edit: extends the Database Schema table
# with 10 new table columns (F:O) describing the additional DB tables, removes
# the old ERD screenshot picture, and tidies up the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database Schema")

# --- 1) Column widths for the newly used columns F:O ---
$ws.Columns("F:F").ColumnWidth = 24.33203125
$ws.Columns("G:G").ColumnWidth = 21.5546875
$ws.Columns("H:H").ColumnWidth = 24.88671875
$ws.Columns("I:I").ColumnWidth = 26.77734375
$ws.Columns("J:J").ColumnWidth = 22.33203125
$ws.Columns("K:K").ColumnWidth = 21.5546875
$ws.Columns("L:L").ColumnWidth = 19.44140625
$ws.Columns("M:M").ColumnWidth = 18.5546875
$ws.Columns("N:N").ColumnWidth = 21.77734375
$ws.Columns("O:O").ColumnWidth = 20

# --- 2) Cell values for the new columns (F:O) added to the Database Schema table ---
$ws.Range("F1").Value = "Fields"
$ws.Range("G1").Value = "Fields"
$ws.Range("H1").Value = "Fields"
$ws.Range("I1").Value = "Fields"
$ws.Range("J1").Value = "Fields"
$ws.Range("K1").Value = "Fields"
$ws.Range("L1").Value = "Fields"
$ws.Range("M1").Value = "Fields"
$ws.Range("N1").Value = "Fields"
$ws.Range("O1").Value = "Fields"
$ws.Range("F2").Value = "food_details"
$ws.Range("G2").Value = "food_inventory"
$ws.Range("H2").Value = "food_requests"
$ws.Range("I2").Value = "deliveries"
$ws.Range("J2").Value = "donations"
$ws.Range("K2").Value = "donation_items"
$ws.Range("L2").Value = "donation_payments"
$ws.Range("M2").Value = "feeds"
$ws.Range("N2").Value = "payments"
$ws.Range("O2").Value = "status_logs"
$ws.Range("F3").Value = "food_id"
$ws.Range("G3").Value = "inventory_id"
$ws.Range("H3").Value = "request_id"
$ws.Range("I3").Value = "delivery_id"
$ws.Range("J3").Value = "donation_id"
$ws.Range("K3").Value = "donation_item_id"
$ws.Range("L3").Value = "donation_payment_id"
$ws.Range("M3").Value = "feed_id"
$ws.Range("N3").Value = "payment_id"
$ws.Range("O3").Value = "log_id"
$ws.Range("F4").Value = "food_name"
$ws.Range("G4").Value = "kitchen_id"
$ws.Range("H4").Value = "user_id"
$ws.Range("I4").Value = "request_id"
$ws.Range("J4").Value = "user_id"
$ws.Range("K4").Value = "donation_id"
$ws.Range("L4").Value = "donation_id"
$ws.Range("M4").Value = "user_id"
$ws.Range("N4").Value = "user_id"
$ws.Range("O4").Value = "payment_id"
$ws.Range("F5").Value = "food_desc"
$ws.Range("G5").Value = "food_id"
$ws.Range("H5").Value = "kitchen_id"
$ws.Range("I5").Value = "donation_id"
$ws.Range("J5").Value = "kitchen_id"
$ws.Range("K5").Value = "food_id"
$ws.Range("L5").Value = "amount"
$ws.Range("M5").Value = "post_type"
$ws.Range("N5").Value = "request_id"
$ws.Range("O5").Value = "old_status"
$ws.Range("F6").Value = "price"
$ws.Range("G6").Value = "quantity"
$ws.Range("H6").Value = "food_id"
$ws.Range("I6").Value = "driver_id"
$ws.Range("J6").Value = "admin_id"
$ws.Range("K6").Value = "quantity"
$ws.Range("L6").Value = "payment_method"
$ws.Range("M6").Value = "content"
$ws.Range("N6").Value = "delivery_id"
$ws.Range("O6").Value = "new_status"
$ws.Range("F7").Value = "Status"
$ws.Range("G7").Value = "unit"
$ws.Range("H7").Value = "request_date"
$ws.Range("I7").Value = "pickup_location"
$ws.Range("J7").Value = "date_donated"
$ws.Range("L7").Value = "reference_number"
$ws.Range("M7").Value = "date_posted"
$ws.Range("N7").Value = "amount"
$ws.Range("O7").Value = "notes"
$ws.Range("G8").Value = "expirey_date"
$ws.Range("H8").Value = "request_type"
$ws.Range("I8").Value = "dropoff_location"
$ws.Range("J8").Value = "donated_location"
$ws.Range("L8").Value = "delivery_id"
$ws.Range("M8").Value = "status"
$ws.Range("N8").Value = "payment_method"
$ws.Range("O8").Value = "change_date"
$ws.Range("G9").Value = "last_updated"
$ws.Range("H9").Value = "status"
$ws.Range("I9").Value = "status"
$ws.Range("N9").Value = "transaction_type"
$ws.Range("O9").Value = "deleted_date"
$ws.Range("H10").Value = "approved_by"
$ws.Range("I10").Value = "delivery_time"
$ws.Range("N10").Value = "payment_status"
$ws.Range("N11").Value = "transaction_date"

# --- 3) Apply matching cell fill styles by copying formats from existing template cells ---
# xlPasteFormats = -4122
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null
$ws.Range("O6").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null
$ws.Range("J7").PasteSpecial(-4122) | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null
$ws.Range("O7").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("M8").PasteSpecial(-4122) | Out-Null
$ws.Range("N8").PasteSpecial(-4122) | Out-Null
$ws.Range("O8").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("N9").PasteSpecial(-4122) | Out-Null
$ws.Range("O9").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null
$ws.Range("N10").PasteSpecial(-4122) | Out-Null
$ws.Range("N11").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null

$ws.Range("C4").Copy()
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("J4").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("L8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 4) Remove the old ERD screenshot picture from the sheet ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# --- 5) Sheet view: zoom to 115% and move the selection to L8 ---
$ws.Activate()
$ws.Range("L8").Select()
$excel.ActiveWindow.Zoom = 115

